# Update cryptos list values (price + 1h volume change) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.885.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "'2.627.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'596.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'153.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "'2.627.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +10.19%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "'27.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "'3.103.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'67.827.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "'2.618.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'11.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").Value = "'370.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "'71.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.78%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").Value = "'1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'574.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").Value = "'1.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").Value = "'158.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("E44").Value = "  +15.45%  "
$ws.Range("E45").Value = "  +6.05%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'40.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").Value = "'155.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'21.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("E51").Value = "  +0.14%  "
